# Regenerate the "K" column (column G) values for the save_data sheet.
# The K column replaces the old "Strike#" derived values with newly
# computed strike counts (std/mean recalculated upstream, s_vals written here).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sVals = @(1,1,0,4,1,5,0,0,2,3,1,1,1,2,1,1,1,1,1,11,2,0,4,2,5,2,5,3,3,7,3,2,2,5,3,4,5,7,5,3,1)

$startRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
